$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the unnecessary "id_scenario" column (column A), shifting remaining
# columns left.
$ws.Columns.Item(1).Delete()
